$wb = $excel.ActiveWorkbook

# "About" sheet: add note to B9 referencing the EU calibration note
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("B9").Value = "For the EU, roughly calibrated against EU scenario from BNEF."

# "GBCGpUNR" sheet: change the calibrated value from 2000 to 400
$wsData = $wb.Worksheets.Item("GBCGpUNR")
$wsData.Range("B2").Value = 400
